# Kate assigned to book database task
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SprintTracking")

# Rename "Sarah" to "Kate" everywhere it is currently the assignee (E2),
# and clear the assignment from E3 (the "Create API..." task).
$ws.Range("E2").Value = "Kate"
$ws.Range("E3").ClearContents()
